$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (old row 22) since the new data has only 20 data rows (rows 2-21)
$ws.Rows.Item(22).Delete()

# Overwrite rows 2-21 with the updated dataset values (numeric columns C-H).
# Column A (timestamp) and B (label) are unchanged by this edit.
$ws.Range("C2").Value = 1.809551620483398
$ws.Range("D2").Value = -2.747114396095276
$ws.Range("E2").Value = -2.346350741386413
$ws.Range("F2").Value = -0.1327104717493057
$ws.Range("G2").Value = 0.0224492978304624
$ws.Range("H2").Value = 0.07544185966253281

$ws.Range("C3").Value = 1.535980415344238
$ws.Range("D3").Value = -2.904552030563355
$ws.Range("E3").Value = -2.271616220474243
$ws.Range("F3").Value = -0.0704022198915481
$ws.Range("G3").Value = 0.0503963828086853
$ws.Range("H3").Value = -0.0980438739061355

$ws.Range("C4").Value = 1.500288486480713
$ws.Range("D4").Value = -2.749492883682251
$ws.Range("E4").Value = -1.919559156894683
$ws.Range("F4").Value = 0.0665843114256858
$ws.Range("G4").Value = 0.0959058403968811
$ws.Range("H4").Value = -0.0583376325666904

$ws.Range("C5").Value = 1.063152790069579
$ws.Range("D5").Value = -2.663525581359863
$ws.Range("E5").Value = -1.874630331993103
$ws.Range("F5").Value = 0.0452040284872055
$ws.Range("G5").Value = 0.1351539343595504
$ws.Range("H5").Value = -0.1539380401372909

$ws.Range("C6").Value = 1.008758783340455
$ws.Range("D6").Value = -2.917640089988708
$ws.Range("E6").Value = -2.169865667819977
$ws.Range("F6").Value = 0.2449569702148437
$ws.Range("G6").Value = 0.4401284158229828
$ws.Range("H6").Value = -0.2344195395708084

$ws.Range("C7").Value = 0.8518548965454079
$ws.Range("D7").Value = -3.11153244972229
$ws.Range("E7").Value = -1.974024677276608
$ws.Range("F7").Value = 0.5047274231910706
$ws.Range("G7").Value = 1.391398310661316
$ws.Range("H7").Value = 1.06214189529419

$ws.Range("C8").Value = -0.7313633918762299
$ws.Range("D8").Value = -2.766938614845273
$ws.Range("E8").Value = -1.49925755262375
$ws.Range("F8").Value = 0.0247400421649217
$ws.Range("G8").Value = -0.4867068827152252
$ws.Range("H8").Value = -1.290605545043945

$ws.Range("C9").Value = -2.322844505310054
$ws.Range("D9").Value = -2.311514568328858
$ws.Range("E9").Value = -1.98325538635254
$ws.Range("F9").Value = 0.8316930532455444
$ws.Range("G9").Value = -1.329548239707947
$ws.Range("H9").Value = 0.3587306141853332

$ws.Range("C10").Value = 0.3075991153717006
$ws.Range("D10").Value = -2.543269753456111
$ws.Range("E10").Value = -2.584808015823368
$ws.Range("F10").Value = 0.7304421067237854
$ws.Range("G10").Value = -2.097558498382568
$ws.Range("H10").Value = -0.8232936263084412

$ws.Range("C11").Value = -1.671347141265869
$ws.Range("D11").Value = 0.0999624729156494
$ws.Range("E11").Value = -4.535521984100342
$ws.Range("F11").Value = 0.2547308206558227
$ws.Range("G11").Value = -0.4392121136188507
$ws.Range("H11").Value = -0.3918700516223907

$ws.Range("C12").Value = 2.258466720581057
$ws.Range("D12").Value = -3.150174045562746
$ws.Range("E12").Value = -2.508739709854125
$ws.Range("F12").Value = 0.09239336848258969
$ws.Range("G12").Value = -0.5294674634933472
$ws.Range("H12").Value = 0.3912591934204101

$ws.Range("C13").Value = -4.284438991546638
$ws.Range("D13").Value = -5.87428689002991
$ws.Range("E13").Value = 0.2990560531616238
$ws.Range("F13").Value = -1.268614411354065
$ws.Range("G13").Value = -1.000597238540649
$ws.Range("H13").Value = -1.631621122360229

$ws.Range("C14").Value = 3.352084398269707
$ws.Range("D14").Value = -4.846620321273795
$ws.Range("E14").Value = -3.800319671630887
$ws.Range("F14").Value = 1.631926536560059
$ws.Range("G14").Value = 8.462469100952148
$ws.Range("H14").Value = -0.6192646622657776

$ws.Range("C15").Value = 2.565014839172318
$ws.Range("D15").Value = -4.19008378982544
$ws.Range("E15").Value = -4.534866142272937
$ws.Range("F15").Value = -1.297935962677002
$ws.Range("G15").Value = -0.3597996234893799
$ws.Range("H15").Value = 0.7938193678855896

$ws.Range("C16").Value = -0.5449800491333003
$ws.Range("D16").Value = -2.96973985433578
$ws.Range("E16").Value = -1.820678830146782
$ws.Range("F16").Value = -0.845132052898407
$ws.Range("G16").Value = 0.3094032406806946
$ws.Range("H16").Value = -0.9367618560791016

$ws.Range("C17").Value = 0.4441701889038138
$ws.Range("D17").Value = -1.205912351608273
$ws.Range("E17").Value = -2.574502897262589
$ws.Range("F17").Value = -1.232115149497986
$ws.Range("G17").Value = 0.111024759709835
$ws.Range("H17").Value = -0.7185302376747131

$ws.Range("C18").Value = 0.6093713760375863
$ws.Range("D18").Value = -1.489246553182612
$ws.Range("E18").Value = -5.765595197677606
$ws.Range("F18").Value = -0.2264782935380935
$ws.Range("G18").Value = -0.9043859839439392
$ws.Range("H18").Value = 1.454012036323547

$ws.Range("C19").Value = -2.423646736145025
$ws.Range("D19").Value = -3.778019905090329
$ws.Range("E19").Value = -3.443504238128652
$ws.Range("F19").Value = 0.5639813542366028
$ws.Range("G19").Value = 0.6252205967903137
$ws.Range("H19").Value = -0.7906123399734497

$ws.Range("C20").Value = -3.03008975982664
$ws.Range("D20").Value = -3.059792947769166
$ws.Range("E20").Value = -0.8526946783065815
$ws.Range("F20").Value = -0.1081231459975242
$ws.Range("G20").Value = 0.0042760567739605
$ws.Range("H20").Value = 1.263422012329102

$ws.Range("C21").Value = 1.929839134216309
$ws.Range("D21").Value = -3.365105152130127
$ws.Range("E21").Value = -1.33200478553772
$ws.Range("F21").Value = 0.1950187236070633
$ws.Range("G21").Value = -0.3927863538265228
$ws.Range("H21").Value = 0.2987131178379059

